$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 249.5
$ws2.Range("O10").Value = 3876
$ws2.Range("N12").Value = 3894.9
$ws2.Range("O12").Value = 235.265

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.01 = 15535.07 pesos`n✅ 15535.07 pesos = 3.99 = 938.37 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText
